$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "kaxh"
$ws.Range("B2").Value = 3114

$ws.Range("B3").Select()
